$d = $word.ActiveDocument

$full       = "Partager avec la communauté autour ce que vous faites régulièrement cela permet aux gens autour de vous de connaître ce que vous faites !"
$firstText  = "Partager avec la communauté autour ce que vous faites régulièrement "
$secondText = "cela permet aux gens autour de vous de connaître ce que vous faites !"

$rng = $d.Content
$rng.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng.Find.Found) {
    $p  = $rng.Paragraphs(1)
    $pr = $p.Range

    # Pull the paragraph's own opening <w:p ...> tag (with its w14:paraId / rsid*
    # attributes) straight out of the live document via the read-only
    # WordOpenXML round-trip, so the rebuilt paragraph keeps its identity
    # instead of relying on hard-coded GUIDs.
    $openXml = $pr.WordOpenXML
    $null = $openXml -match '<w:p\s[^>]*>'
    $pOpenTag = $matches[0]

    $newInner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">' + $firstText + '</w:t></w:r>' + `
        '<w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>' + $secondText + '</w:t></w:r>'

    $xml = $pOpenTag + $newInner + '</w:p>'

    $pr.InsertXML($xml) | Out-Null
}
